$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1) "Faire un bâtiment simple " bullet: strike the whole paragraph through.
# --------------------------------------------------------------------------
$rBatiment = $d.Content
$null = $rBatiment.Find.Execute("Faire un bâtiment simple", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)
$pBatiment = $rBatiment.Paragraphs(1)
$pBatiment.Range.Font.StrikeThrough = 1

# --------------------------------------------------------------------------
# 2) "Créer un personnage..." bullet: append a new trailing run of text.
# --------------------------------------------------------------------------
$rPersonnage = $d.Content
$null = $rPersonnage.Find.Execute("Créer un personnage et le faire mouvoir dans cette environnement", `
                                   $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pPersonnage = $rPersonnage.Paragraphs(1)
$endPersonnage = $pPersonnage.Range.End - 1
$addedText = " (j’ai retirer la gravité)"
$insertionPoint = $d.Range($endPersonnage, $endPersonnage)
$insertionPoint.InsertAfter($addedText)
# The freshly inserted text shares the same run formatting as what precedes
# it, so the engine would otherwise fold it straight back into that run.
# Nudging a character property away and back forces it to persist as its
# own distinct <w:r>.
$rNewRun = $d.Range($endPersonnage, $endPersonnage + $addedText.Length)
$rNewRun.Font.Size = 12
$rNewRun.Font.Size = 18

# --------------------------------------------------------------------------
# 3) "Mettre un Timer..." bullet: drop the stale spell-check run-splits and
#    leave just two clean runs (matching the two sentences/clauses).
# --------------------------------------------------------------------------
$rTimer = $d.Content
$null = $rTimer.Find.Execute("Timer déclenchable", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$pTimer = $rTimer.Paragraphs(1)
$timerStart = $pTimer.Range.Start
$timerTextLen = $pTimer.Range.Text.Length - 1
$rTimerWhole = $d.Range($timerStart, $timerStart + $timerTextLen)
$timerWholeText = $rTimerWhole.Text
$rTimerWhole.Find.Execute($timerWholeText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $timerWholeText, 2)

$firstHalf = "Mettre un Timer déclenchable à l’aide d’un touche"
$splitPoint = $timerStart + $firstHalf.Length
$rSecondHalf = $d.Range($splitPoint, $timerStart + $timerTextLen)
$rSecondHalf.Font.Size = 12
$rSecondHalf.Font.Size = 18

# --------------------------------------------------------------------------
# 4) "Rendre la simulation..." bullet: collapse everything into one run.
# --------------------------------------------------------------------------
$rRendre = $d.Content
$null = $rRendre.Find.Execute("Rendre la simulation", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
$pRendre = $rRendre.Paragraphs(1)
$rendreStart = $pRendre.Range.Start
$rendreTextLen = $pRendre.Range.Text.Length - 1
$rRendreWhole = $d.Range($rendreStart, $rendreStart + $rendreTextLen)
$rendreWholeText = $rRendreWhole.Text
$rRendreWhole.Find.Execute($rendreWholeText, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $rendreWholeText, 2)
